$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates as serial numbers, matching existing column A style)
$newRows = @(
    @{ Row = 227; A = 44301; B = 5;  C = 38; D = 157.53254290689 },
    @{ Row = 228; A = 44302; B = 10; C = 42; D = 174.1149158444573 },
    @{ Row = 229; A = 44303; B = 12; C = 46; D = 190.6972887820247 }
)

# Template cell used to clone the date-column formatting (style index 2 in the original file)
$templateA = $ws.Range("A226")

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $cellA = $ws.Range("A" + $rowNum)
    $templateA.Copy()
    $cellA.PasteSpecial(-4122)  # xlPasteFormats
    $cellA.Value = $r.A

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
}
